$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Add($ws.Range("B24"), "https://example.com/standoff", "", "Standoff", "Standoff")
Write-Output "Count all after add: $($ws.Hyperlinks.Count)"
